# Adds two new pharmacy-item rows (FEROGLOBIN 30 CAPS and LEZBERG TRIO
# 20/5/12.5 TAB) above the totals/footer rows, shifting the totals row
# and the footer row down, and updates the grand total accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 is the first (and, before this edit, only) item row.
# Row 5 is the totals row and row 6 is the footer row.
# Insert two fresh rows right before the totals row so the new items can
# be added as rows 5 and 6; the totals/footer rows shift down to 7/8.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(6).Insert()

# New rows come in blank/unformatted; copy the formatting (styles +
# column layout) from the existing item row (row 4) down onto them.
$ws.Range("A4:N4").Copy()
$ws.Range("A5:N5").PasteSpecial(-4122)
$ws.Range("A4:N4").Copy()
$ws.Range("A6:N6").PasteSpecial(-4122)

# Recreate the merged cell layout for the two new item rows (matches the
# merge layout already used by row 4).
$ws.Range("B5:G5").Merge()
$ws.Range("H5:K5").Merge()
$ws.Range("L5:M5").Merge()
$ws.Range("B6:G6").Merge()
$ws.Range("H6:K6").Merge()
$ws.Range("L6:M6").Merge()

# Restore the correct row heights for the new item rows.
$ws.Rows.Item(5).RowHeight = 25.5
$ws.Rows.Item(6).RowHeight = 24.75

# Fill in the data for the new item rows.
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "FEROGLOBIN 30 CAPS"
$ws.Range("H5").Value = "0:1"
$ws.Range("L5").Value = 90
$ws.Range("N5").Value = "0:2"

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "LEZBERG TRIO 20/5/12.5 TAB"
$ws.Range("H6").Value = "0:2"
$ws.Range("L6").Value = 38
$ws.Range("N6").Value = "0:0"

# Update the grand total (was 62, now 62 + 90 + 38 = 190) which now lives
# on row 7 after the insert.
$ws.Range("K7").Value = 190

Write-Host "edit applied"
